$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 196: Gavin Stone
$ws.Range("A196").Value = "Gavin Stone"
$ws.Range("B196").Value = "https://midfield.mlbstatic.com/v1/people/694813/spots/120"
$ws.Range("C196").Value = "RHP"
$ws.Hyperlinks.Add($ws.Range("B196"), "https://midfield.mlbstatic.com/v1/people/694813/spots/120")
$ws.Range("B196").Style = "Hyperlink"

# Row 197: Luis Medina
$ws.Range("A197").Value = "Luis Medina"
$ws.Range("B197").Value = "https://midfield.mlbstatic.com/v1/people/665622/spots/120"
$ws.Range("C197").Value = "RHP"
$ws.Hyperlinks.Add($ws.Range("B197"), "https://midfield.mlbstatic.com/v1/people/665622/spots/120")
$ws.Range("B197").Style = "Hyperlink"

# Row 198: Luis Ortiz (no handedness listed)
$ws.Range("A198").Value = "Luis Ortiz"
$ws.Range("B198").Value = "https://midfield.mlbstatic.com/v1/people/123456/spots/120"
$ws.Hyperlinks.Add($ws.Range("B198"), "https://midfield.mlbstatic.com/v1/people/123456/spots/120")
$ws.Range("B198").Style = "Hyperlink"

# Update view state to match the end of the newly added data
$ws.Range("C198").Select()
$excel.ActiveWindow.ScrollRow = 173
